# Append new incident rows (186-196) to the management log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("WC47 NACP", "Fallo tolva",                    "2024-06-10", "12:19:24", "Mañana", "12:19:27", "0:00:03", "-0.01 minutos"),
    @("WC48 P5F",  "AOI (malla)",                     "2024-06-10", "12:24:10", "Mañana", "12:24:12", "0:00:02", "-0.00 minutos"),
    @("WC48 P5F",  "AOI (fallo etiqueta)",             "2024-06-10", "12:25:47", "Mañana", "12:25:48", "0:00:01", "0.40 minutos"),
    @("WC48 P5F",  "AOI (fallo etiqueta)",             "2024-06-10", "12:26:00", "Mañana", "12:26:00", "0:00:00", "0.30 minutos"),
    @("WC48 P5F",  "Etiquetadora",                     "2024-06-10", "12:30:37", "Mañana", "12:30:38", "0:00:01", "-0.00 minutos"),
    @("WC47 NACP", "Etiquetadora",                     "2024-06-10", "12:38:44", "Mañana", "12:38:48", "0:00:04", "-0.01 minutos"),
    @("WC47 NACP", "Ascensor no sube",                 "2024-06-10", "12:42:58", "Mañana", "12:42:59", "0:00:01", "-0.00 minutos"),
    @("WC49 P5H",  "No coloca bien el sealling",       "2024-06-10", "12:50:25", "Mañana", "12:50:26", "0:00:01", "-0.00 minutos"),
    @("WC49 P5H",  "La cámara no detecta Busbar",      "2024-06-10", "12:50:31", "Mañana", "12:50:33", "0:00:02", "0.02 minutos"),
    @("WC49 P5H",  "La cámara no detecta Busbar",      "2024-06-10", "12:50:32", "Mañana", "12:50:33", "0:00:01", "0.02 minutos"),
    @("WC47 NACP", "No atornilla clips",               "2024-06-10", "12:54:48", "Mañana", "12:54:50", "0:00:02", "-0.01 minutos")
)

$startRow = 186
$endRow = $startRow + $newRows.Count - 1

# Force the new range to be stored as plain text so date/time-looking
# strings (e.g. "2024-06-10", "12:19:24") are not auto-converted into
# date/time serial values by Excel.
$fillRange = $ws.Range("A" + $startRow + ":H" + $endRow)
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
